# Generate Report for Handback
# Adds the handback-attempt result for the "92930ec5-3eb3-496a-b82d-ef06c97dbc52"
# row (row 6) on both the "zh-cn" and "de-de" sheets: a hyperlinked "Latest Target
# File" (I6), the handback xlf name (J6), the handback datetime (K6) and an
# "Error Detail" (P6) explaining the handback file was stale. Also widens the
# "Error Detail" column (P) so the long message is readable.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/acc69ae610172f165f58d58a450285237dae3442/e2e/92930ec5-3eb3-496a-b82d-ef06c97dbc52.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0a635a6125c6decf100a2a09e0d4b4c88edafb4b/e2e/92930ec5-3eb3-496a-b82d-ef06c97dbc52.md."
$currentUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/acc69ae610172f165f58d58a450285237dae3442/e2e/92930ec5-3eb3-496a-b82d-ef06c97dbc52.md"
$fileDisplay = "92930ec5-3eb3-496a-b82d-ef06c97dbc52.md"

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Columns.Item(16).ColumnWidth = 39.14

$ws.Hyperlinks.Add($ws.Range("I6"), $currentUrl, "", "", $fileDisplay)
$ws.Range("J6").Value = "92930ec5-3eb3-496a-b82d-ef06c97dbc52.17c664f3e183ac20bdf299bcf6755d8f8db774d7.zh-cn.xlf"
$ws.Range("K6").Value = "2016-08-28 04:42:47"
$ws.Range("P6").Value = $errorDetail

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")

$ws.Columns.Item(16).ColumnWidth = 39.14

$ws.Hyperlinks.Add($ws.Range("I6"), $currentUrl, "", "", $fileDisplay)
$ws.Range("J6").Value = "92930ec5-3eb3-496a-b82d-ef06c97dbc52.17c664f3e183ac20bdf299bcf6755d8f8db774d7.de-de.xlf"
$ws.Range("K6").Value = "2016-08-28 04:42:54"
$ws.Range("P6").Value = $errorDetail
